$wb = $excel.ActiveWorkbook

# New trade row data (Trade #20) to append as row 21 on both the
# "All Trades" and "base_strategy" sheets.
$tradeNum   = 20
$tradeDate  = "2026-02-16"
$tradeTime  = "22:54:04"
$strategy   = "base_strategy"
$side       = "UP"
$entryPrice = 49.999998
$exitPrice  = ""
$status     = "OPEN"
$pnlPct     = 0
$pnlDollar  = 0
$capAfter   = 100
$entrySlip  = 0
$exitSlip   = 0
$confidence = 0.6
$entryReason = "Normal spread capture: 19600 bps"
$exitReason  = ""
$duration    = 0

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $r = 21
    $prev = $r - 1

    $ws.Cells.Item($r, 1).Value  = $tradeNum

    # Column B holds a "YYYY-MM-DD" looking string; Excel would otherwise
    # auto-convert it to a date serial number. Force it to stay literal
    # text (matching the source file's string representation), then
    # reset the cell style back to Normal so no extra style is retained.
    $dateCell = $ws.Cells.Item($r, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $tradeDate
    $dateCell.Style = "Normal"

    $ws.Cells.Item($r, 3).Value  = $tradeTime
    $ws.Cells.Item($r, 4).Value  = $strategy
    $ws.Cells.Item($r, 5).Value  = $side
    $ws.Cells.Item($r, 6).Value  = $entryPrice

    # Column G (Exit Price) is blank for an OPEN trade. Assigning an empty
    # string does not materialize the cell, so copy the blank cell above
    # it (which already exists in the sheet) down into this row instead.
    $ws.Range("G" + $prev).Copy($ws.Range("G" + $r))

    $ws.Cells.Item($r, 8).Value  = $status
    $ws.Cells.Item($r, 9).Value  = $pnlPct
    $ws.Cells.Item($r, 10).Value = $pnlDollar
    $ws.Cells.Item($r, 11).Value = $capAfter
    $ws.Cells.Item($r, 12).Value = $entrySlip
    $ws.Cells.Item($r, 13).Value = $exitSlip
    $ws.Cells.Item($r, 14).Value = $confidence
    $ws.Cells.Item($r, 15).Value = $entryReason

    # Column P (Exit Reason) is likewise blank; reuse the same technique.
    $ws.Range("P" + $prev).Copy($ws.Range("P" + $r))

    $ws.Cells.Item($r, 17).Value = $duration
}
